# Auto-generated Excel COM-interop script
# Updates market-price-derived columns (H-N) for specific Leve rows across sheets,
# reflecting a scheduled market-data refresh (no formulas involved; plain values).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 218.5625
$ws.Range("I9").Value = 317
$ws.Range("J9").Value = 92
$ws.Range("K9").Value = 317
$ws.Range("L9").Value = 92
$ws.Range("M9").Value = -148
$ws.Range("N9").Value = -430

$ws.Range("I19").Value = 3416.2104
$ws.Range("J19").Value = 3568.5715
$ws.Range("K19").Value = 3416.2104
$ws.Range("L19").Value = 3568.5715
$ws.Range("M19").Value = -3241.2104
$ws.Range("N19").Value = -3918.5715

$ws.Range("H32").Value = 7435.4346
$ws.Range("J32").Value = 4926.385
$ws.Range("L32").Value = 4926.385
$ws.Range("N32").Value = -5578.385

$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3825

$ws.Range("H43").Value = 10007.1
$ws.Range("J43").Value = 5008.875
$ws.Range("L43").Value = 5008.875
$ws.Range("N43").Value = -5146.875

$ws.Range("H80").Value = 6570.0586
$ws.Range("I80").Value = 7391.5713
$ws.Range("J80").Value = 5995
$ws.Range("K80").Value = 22174.7139
$ws.Range("L80").Value = 17985
$ws.Range("M80").Value = -21176.7139
$ws.Range("N80").Value = -19981

$ws.Range("H83").Value = 6570.0586
$ws.Range("I83").Value = 7391.5713
$ws.Range("J83").Value = 5995
$ws.Range("K83").Value = 66524.14169999999
$ws.Range("L83").Value = 53955
$ws.Range("M83").Value = -61532.14169999999
$ws.Range("N83").Value = -63939

$ws.Range("H92").Value = 564.8333
$ws.Range("I92").Value = 446.1875
$ws.Range("J92").Value = 1514
$ws.Range("K92").Value = 446.1875
$ws.Range("L92").Value = 1514
$ws.Range("M92").Value = 801.8125
$ws.Range("N92").Value = -4010

$ws.Range("H112").Value = 8566.212
$ws.Range("J112").Value = 8566.212
$ws.Range("L112").Value = 25698.636
$ws.Range("N112").Value = -27914.636

$ws.Range("H116").Value = 5621.1333
$ws.Range("I116").Value = 5031.3335
$ws.Range("J116").Value = 6014.3335
$ws.Range("K116").Value = 5031.3335
$ws.Range("L116").Value = 6014.3335
$ws.Range("M116").Value = -1589.3335
$ws.Range("N116").Value = -12898.3335

$ws.Range("H138").Value = 2763.4866
$ws.Range("I138").Value = 1430.5
$ws.Range("J138").Value = 4026.3157
$ws.Range("K138").Value = 4291.5
$ws.Range("L138").Value = 12078.9471
$ws.Range("M138").Value = 848.5
$ws.Range("N138").Value = -22358.9471

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22249.896
$ws.Range("I32").Value = 22419.553
$ws.Range("J32").Value = 17499.5
$ws.Range("K32").Value = 22419.553
$ws.Range("L32").Value = 17499.5
$ws.Range("M32").Value = -22132.553
$ws.Range("N32").Value = -18073.5

$ws.Range("H45").Value = 4266.4
$ws.Range("I45").Value = 2381
$ws.Range("K45").Value = 2381
$ws.Range("M45").Value = -2004

$ws.Range("H74").Value = 547920.2
$ws.Range("I74").Value = 1000770.7
$ws.Range("J74").Value = 4499.6
$ws.Range("K74").Value = 1000770.7
$ws.Range("L74").Value = 4499.6
$ws.Range("M74").Value = -999896.7
$ws.Range("N74").Value = -6247.6

$ws.Range("H77").Value = 547920.2
$ws.Range("I77").Value = 1000770.7
$ws.Range("J77").Value = 4499.6
$ws.Range("K77").Value = 5003853.5
$ws.Range("L77").Value = 22498
$ws.Range("M77").Value = -4999485.5
$ws.Range("N77").Value = -31234

$ws.Range("H132").Value = 3256.476
$ws.Range("I132").Value = 1270.5454
$ws.Range("J132").Value = 3961.1614
$ws.Range("K132").Value = 3811.6362
$ws.Range("L132").Value = 11883.4842
$ws.Range("M132").Value = -1281.6362
$ws.Range("N132").Value = -16943.4842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 496.35715
$ws.Range("I22").Value = 503.76923
$ws.Range("K22").Value = 503.76923
$ws.Range("M22").Value = -330.76923

$ws.Range("H94").Value = 508.3846
$ws.Range("I94").Value = 527
$ws.Range("J94").Value = 446.33334
$ws.Range("K94").Value = 527
$ws.Range("L94").Value = 446.33334
$ws.Range("M94").Value = -76
$ws.Range("N94").Value = -1348.33334

$ws.Range("H99").Value = 3028.1428
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 4665.6665
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 4665.6665
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -7661.6665

$ws.Range("H107").Value = 65050
$ws.Range("I107").Value = 85983.336
$ws.Range("J107").Value = 2250
$ws.Range("K107").Value = 85983.336
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = -84063.336
$ws.Range("N107").Value = -6090

$ws.Range("H140").Value = 103662.13
$ws.Range("J140").Value = 103662.13
$ws.Range("L140").Value = 103662.13
$ws.Range("N140").Value = -114022.13

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1293.6666
$ws.Range("I16").Value = 1080.75
$ws.Range("K16").Value = 1080.75
$ws.Range("M16").Value = -793.75

$ws.Range("H31").Value = 7144614
$ws.Range("I31").Value = 8335049.5
$ws.Range("J31").Value = 1999.5
$ws.Range("K31").Value = 8335049.5
$ws.Range("L31").Value = 1999.5
$ws.Range("M31").Value = -8334754.5
$ws.Range("N31").Value = -2589.5

$ws.Range("H34").Value = 7144614
$ws.Range("I34").Value = 8335049.5
$ws.Range("J34").Value = 1999.5
$ws.Range("K34").Value = 8335049.5
$ws.Range("L34").Value = 1999.5
$ws.Range("M34").Value = -8334847.5
$ws.Range("N34").Value = -2403.5

$ws.Range("H99").Value = 6740.08
$ws.Range("I99").Value = 4587.6665
$ws.Range("J99").Value = 9968.700000000001
$ws.Range("K99").Value = 4587.6665
$ws.Range("L99").Value = 9968.700000000001
$ws.Range("M99").Value = -3089.6665
$ws.Range("N99").Value = -12964.7

$ws.Range("H113").Value = 1293.6666
$ws.Range("I113").Value = 1080.75
$ws.Range("K113").Value = 1080.75
$ws.Range("M113").Value = 1089.25

$ws.Range("H126").Value = 6740.08
$ws.Range("I126").Value = 4587.6665
$ws.Range("J126").Value = 9968.700000000001
$ws.Range("K126").Value = 13762.9995
$ws.Range("L126").Value = 29906.1
$ws.Range("M126").Value = -11292.9995
$ws.Range("N126").Value = -34846.10000000001

$ws.Range("H134").Value = 3320.1428
$ws.Range("I134").Value = 2965.7
$ws.Range("K134").Value = 8897.099999999999
$ws.Range("M134").Value = -6362.099999999999

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 184.6
$ws.Range("I14").Value = 184.6
$ws.Range("K14").Value = 553.8
$ws.Range("M14").Value = -380.8

$ws.Range("H131").Value = 1965
$ws.Range("J131").Value = 2319
$ws.Range("L131").Value = 6957
$ws.Range("N131").Value = -17037

$ws.Range("H132").Value = 1555.9286
$ws.Range("I132").Value = 1589.2727
$ws.Range("J132").Value = 1433.6666
$ws.Range("K132").Value = 14303.4543
$ws.Range("L132").Value = 12902.9994
$ws.Range("N132").Value = -17962.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 12499.5
$ws.Range("I35").Value = 10000
$ws.Range("J35").Value = 14999
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 14999
$ws.Range("M35").Value = -9702
$ws.Range("N35").Value = -15595

$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("N42").Value = -50970

$ws.Range("H102").Value = 1471.0769
$ws.Range("I102").Value = 1394.8223
$ws.Range("K102").Value = 1394.8223
$ws.Range("M102").Value = 227.1777

$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("N115").Value = -52350

$ws.Range("H132").Value = 3735.6333
$ws.Range("I132").Value = 3310.7693
$ws.Range("K132").Value = 9932.3079
$ws.Range("M132").Value = -7402.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1908.9166
$ws.Range("I22").Value = 1857.8572
$ws.Range("K22").Value = 1857.8572
$ws.Range("M22").Value = -1562.8572

$ws.Range("H27").Value = 1908.9166
$ws.Range("I27").Value = 1857.8572
$ws.Range("K27").Value = 1857.8572
$ws.Range("M27").Value = -1750.8572

$ws.Range("H46").Value = 5385.2354
$ws.Range("I46").Value = 1160.6
$ws.Range("K46").Value = 1160.6
$ws.Range("M46").Value = -972.5999999999999

$ws.Range("H61").Value = 981.2222
$ws.Range("I61").Value = 981.2222
$ws.Range("K61").Value = 981.2222
$ws.Range("M61").Value = -779.2222

$ws.Range("H82").Value = 4808.8335
$ws.Range("I82").Value = 3541.4375
$ws.Range("J82").Value = 7343.625
$ws.Range("K82").Value = 3541.4375
$ws.Range("L82").Value = 7343.625
$ws.Range("M82").Value = -3180.4375
$ws.Range("N82").Value = -8065.625

$ws.Range("H85").Value = 4808.8335
$ws.Range("I85").Value = 3541.4375
$ws.Range("J85").Value = 7343.625
$ws.Range("K85").Value = 3541.4375
$ws.Range("L85").Value = 7343.625
$ws.Range("M85").Value = -2293.4375
$ws.Range("N85").Value = -9839.625

$ws.Range("H113").Value = 981.2222
$ws.Range("I113").Value = 981.2222
$ws.Range("K113").Value = 981.2222
$ws.Range("M113").Value = 1188.7778

$ws.Range("H132").Value = 4858.7334
$ws.Range("I132").Value = 4991.5713
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 14974.7139
$ws.Range("L132").Value = 8997
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 702.5333000000001
$ws.Range("I107").Value = 729.0769
$ws.Range("K107").Value = 2187.2307
$ws.Range("M107").Value = -267.2307000000001

$ws.Range("H126").Value = 4854.4287
$ws.Range("I126").Value = 4996.8335
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 14990.5005
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 4740.706
$ws.Range("I132").Value = 5114.4287
$ws.Range("K132").Value = 15343.2861
$ws.Range("M132").Value = -12813.2861
